$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.770.44"
$ws.Range("E2").Value = "  +0.66%  "

$ws.Range("D3").Value = "3.804.23"
$ws.Range("E3").Value = "  +0.95%  "

$ws.Range("E4").Value = "  -0.01%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "597.85"
$ws.Range("E5").Value = "  +0.78%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "167.50"
$ws.Range("E6").Value = "  +0.87%  "

$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.521"
$ws.Range("E8").Value = "  -0.08%  "

$ws.Range("B9").Value = "Dogecoin"
$ws.Range("C9").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.160"
$ws.Range("E9").Value = "  +1.05%  "

$ws.Range("B10").Value = "Toncoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.30"
$ws.Range("E10").Value = "  -1.57%  "

$ws.Range("B11").Value = "Cardano"
$ws.Range("C11").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.450"
$ws.Range("E11").Value = "  -0.22%  "

$ws.Range("B12").Value = "ShibaInu"
$ws.Range("C12").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.0000254"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("B13").Value = "Avalanche"
$ws.Range("C13").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "36.06"
$ws.Range("E13").Value = "  +0.57%  "

$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.441.77"
$ws.Range("E14").Value = "  +0.76%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "3.798.09"
$ws.Range("E15").Value = "  +0.60%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "18.61"
$ws.Range("E16").Value = "  +5.09%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "67.783.50"
$ws.Range("E17").Value = "  +0.48%  "

$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "7.11"
$ws.Range("E18").Value = "  +2.71%  "

$ws.Range("B19").Value = "TRON"
$ws.Range("C19").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.112"
$ws.Range("E19").Value = "  +0.19%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "461.76"
$ws.Range("E20").Value = "  +0.83%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.94"
$ws.Range("E21").Value = "  -6.02%  "

$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.702"
$ws.Range("E22").Value = "  +0.95%  "

$ws.Range("B23").Value = "PEPE"
$ws.Range("C23").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.0000154"
$ws.Range("E23").Value = "  +0.56%  "

$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "83.59"
$ws.Range("E24").Value = "  +0.43%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "12.13"
$ws.Range("E25").Value = "  +2.57%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.11"
$ws.Range("E26").Value = "  -1.42%  "

$ws.Range("B27").Value = "RenderToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.03"
$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E28").Value = "  +0.08%  "

$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.953.60"
$ws.Range("E29").Value = "  +0.94%  "

$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "2.79"
$ws.Range("E30").Value = "  +1.36%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.26"
$ws.Range("E31").Value = "  +4.53%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "7.29"
$ws.Range("E32").Value = "  +1.22%  "

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "29.71"
$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("B34").Value = "Aptos"
$ws.Range("C34").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "9.12"
$ws.Range("E34").Value = "  -0.53%  "

$ws.Range("B35").Value = "Binance-PegBSC-USD"
$ws.Range("C35").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.996"
$ws.Range("E35").Value = "  -0.39%  "

$ws.Range("B36").Value = "RenzoRestakedETH"
$ws.Range("C36").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D36").Value = "3.744.87"
$ws.Range("E36").Value = "  +0.63%  "

$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.100"
$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "3.43"
$ws.Range("E38").Value = "  +0.75%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.138"
$ws.Range("E39").Value = "  +0.54%  "

$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E40").Value = "  +0.76%  "

$ws.Range("B41").Value = "Filecoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "5.79"
$ws.Range("E41").Value = "  +1.09%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E42").Value = "  -0.01%  "

$ws.Range("B43").Value = "USDe"
$ws.Range("C43").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "1.00"
$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "48.18"
$ws.Range("E44").Value = "  +2.94%  "

$ws.Range("B45").Value = "Arweave"
$ws.Range("C45").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "43.75"
$ws.Range("E45").Value = "  -0.87%  "

$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.299"
$ws.Range("E46").Value = "  -0.07%  "

$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "149.38"
$ws.Range("E47").Value = "  +2.26%  "

$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "400.37"
$ws.Range("E48").Value = "  +2.48%  "

$ws.Range("E49").Value = "  +0.00%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "1.84"
$ws.Range("E50").Value = "  -3.07%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "26.72"
$ws.Range("E51").Value = "  +6.04%  "
